$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "2.0.0-sd-202406-matchbox-patch"
$meta.Range("B8").Value = "2024-06-19T17:47:42+02:00"
$meta.Range("B10").Value = "HL7 International - Structured Documents (http://www.hl7.org/Special/committees/structure, structdog@lists.HL7.org)"

# --- Elements sheet updates ---
$elem = $wb.Worksheets.Item("Elements")
$elem.Range("Z12").Value = "http://hl7.org/cda/stds/core/ValueSet/CDAActClass"
$elem.Range("Z13").Value = "http://hl7.org/cda/stds/core/ValueSet/CDAActMood"
# Min / Base Min are stored as text "0" (shared string), not numeric 0 -
# force text interpretation the same way a user would in Excel (leading apostrophe)
$elem.Range("F12").Value = "'0"
$elem.Range("AG12").Value = "'0"
